$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: drop the bold weight but keep the border/alignment formatting.
$ws.Range("A1:F1").Font.Bold = $false

# id_kho column on the original rows becomes bold (no border).
$ws.Range("A2:A13").ClearFormats()
$ws.Range("A2:A13").Font.Bold = $true

# New warehouse rows appended below the existing 13 rows (rows 14-25).
$newRows = @(
    @(21283000, "Kho Giao Hàng Nặng Quảng Bình"),
    @(21086000, "Kho Giao Hàng Nặng Thanh Hóa"),
    @(21089000, "Kho Giao Hàng Nặng Đà Nẵng"),
    @(21096000, "Kho Giao Hàng Nặng Huế"),
    @(21284000, "Kho Giao Hàng Nặng Quảng Ngãi"),
    @(21090000, "Kho Giao Hàng Nặng Đắk Lắk"),
    @(21525000, "Kho Giao Hàng Nặng Đắk Nông"),
    @(21091000, "Kho Giao Hàng Nặng Gia Lai"),
    @(21087000, "Kho Giao Hàng Nặng Bình Định"),
    @(21285000, "Kho Giao Hàng Nặng Bình Thuận"),
    @(21094000, "Kho Giao Hàng Nặng Khánh Hòa"),
    @(22168000, "Kho Giao Hàng Nặng Hoài Nhơn-Bình Định")
)

$r = 14
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 6).Value = "Ca chiều"
    $r++
}

$ws.Range("H21").Select()
